# Add data for 2025-12-31
# Updates 2025 (column L) violent-crime counts across the Citywide Totals sheet,
# the By Neighborhood summary sheet, and the individual neighborhood sheets that
# received new/updated records for that date.

$wb = $excel.ActiveWorkbook

# --- Citywide Totals ---
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 6565   # was 6556
$ws.Range("L3").Value = 7095   # was 7079
$ws.Range("B4").Value = 1722   # was 1721
$ws.Range("H4").Value = 1768   # was 1767
$ws.Range("J4").Value = 1885   # was 1884
$ws.Range("L4").Value = 1769   # was 1766
$ws.Range("L5").Value = 422   # was 419
$ws.Range("L6").Value = 5811   # was 5795
$ws.Range("B7").Value = 23354   # was 23353
$ws.Range("H7").Value = 26084   # was 26083
$ws.Range("J7").Value = 29364   # was 29363
$ws.Range("L7").Value = 21662   # was 21615

# --- By Neighborhood ---
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L3").Value = 14   # was 13
$ws.Range("L7").Value = 696   # was 693
$ws.Range("L8").Value = 1428   # was 1425
$ws.Range("L9").Value = 124   # was 123
$ws.Range("L19").Value = 598   # was 596
$ws.Range("L20").Value = 544   # was 543
$ws.Range("L23").Value = 229   # was 227
$ws.Range("L29").Value = 1210   # was 1209
$ws.Range("L33").Value = 969   # was 968
$ws.Range("L37").Value = 835   # was 833
$ws.Range("L41").Value = 93   # was 92
$ws.Range("L42").Value = 683   # was 680
$ws.Range("L47").Value = 152   # was 151
$ws.Range("L51").Value = 266   # was 265
$ws.Range("L57").Value = 73   # was 72
$ws.Range("L58").Value = 11   # was 12
$ws.Range("B63").Value = 426   # was 425
$ws.Range("H63").Value = 319   # was 318
$ws.Range("J63").Value = 240   # was 239
$ws.Range("L63").Value = 69   # was 70
$ws.Range("L65").Value = 433   # was 432
$ws.Range("L67").Value = 755   # was 753
$ws.Range("L71").Value = 55   # was 54
$ws.Range("L73").Value = 171   # was 169
$ws.Range("L76").Value = 346   # was 343
$ws.Range("L77").Value = 144   # was 143
$ws.Range("L78").Value = 285   # was 283
$ws.Range("L79").Value = 603   # was 601
$ws.Range("L83").Value = 477   # was 476
$ws.Range("L85").Value = 1080   # was 1079
$ws.Range("L86").Value = 136   # was 135
$ws.Range("L87").Value = 59   # was 57
$ws.Range("L89").Value = 290   # was 289
$ws.Range("L91").Value = 293   # was 291
$ws.Range("L95").Value = 302   # was 301
$ws.Range("L96").Value = 239   # was 238
$ws.Range("L99").Value = 375   # was 374
$ws.Range("B101").Value = 23354   # was 23353
$ws.Range("H101").Value = 26084   # was 26083
$ws.Range("J101").Value = 29364   # was 29363
$ws.Range("L101").Value = 21662   # was 21615

# --- West Ridge ---
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L3").Value = 69   # was 68
$ws.Range("L7").Value = 239   # was 238

# --- Auburn Gresham ---
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L4").Value = 46   # was 45
$ws.Range("L5").Value = 20   # was 19
$ws.Range("L6").Value = 165   # was 164
$ws.Range("L7").Value = 696   # was 693

# --- Uptown ---
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L4").Value = 46   # was 45
$ws.Range("L7").Value = 290   # was 289

# --- South Shore ---
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L5").Value = 23   # was 22
$ws.Range("L7").Value = 1080   # was 1079

# --- Austin ---
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 436   # was 435
$ws.Range("L3").Value = 503   # was 502
$ws.Range("L6").Value = 344   # was 343
$ws.Range("L7").Value = 1428   # was 1425

# --- South Chicago ---
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L3").Value = 189   # was 188
$ws.Range("L7").Value = 477   # was 476

# --- Garfield Park ---
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 265   # was 264
$ws.Range("L7").Value = 969   # was 968

# --- West Pullman ---
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L6").Value = 70   # was 69
$ws.Range("L7").Value = 302   # was 301

# --- Grand Crossing ---
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 252   # was 251
$ws.Range("L6").Value = 215   # was 214
$ws.Range("L7").Value = 835   # was 833

# --- New City ---
$ws = $wb.Worksheets.Item("New City")
$ws.Range("L6").Value = 105   # was 104
$ws.Range("L7").Value = 433   # was 432

# --- Woodlawn ---
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L6").Value = 80   # was 79
$ws.Range("L7").Value = 375   # was 374

# --- North Lawndale ---
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 293   # was 292
$ws.Range("L6").Value = 177   # was 176
$ws.Range("L7").Value = 755   # was 753

# --- Englewood ---
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 369   # was 368
$ws.Range("L3").Value = 468   # was 467
$ws.Range("L4").Value = 65   # was 66
$ws.Range("L7").Value = 1210   # was 1209

# --- Chatham ---
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 213   # was 212
$ws.Range("L3").Value = 181   # was 180
$ws.Range("L7").Value = 598   # was 596

# --- River North ---
$ws = $wb.Worksheets.Item("River North")
$ws.Range("L6").Value = 157   # was 154
$ws.Range("L7").Value = 346   # was 343

# --- Hermosa ---
$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L2").Value = 29   # was 28
$ws.Range("L7").Value = 93   # was 92

# --- Humboldt Park ---
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 185   # was 184
$ws.Range("L4").Value = 52   # was 51
$ws.Range("L5").Value = 18   # was 17
$ws.Range("L7").Value = 683   # was 680

# --- Rogers Park ---
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L2").Value = 75   # was 74
$ws.Range("L3").Value = 93   # was 92
$ws.Range("L7").Value = 285   # was 283

# --- Douglas ---
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L4").Value = 22   # was 21
$ws.Range("L6").Value = 58   # was 57
$ws.Range("L7").Value = 229   # was 227

# --- Washington Park ---
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 132   # was 130
$ws.Range("L7").Value = 293   # was 291

# --- Roseland ---
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 196   # was 195
$ws.Range("L6").Value = 161   # was 160
$ws.Range("L7").Value = 603   # was 601

# --- Chicago Lawn ---
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 191   # was 190
$ws.Range("L7").Value = 544   # was 543

# --- Kenwood ---
$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L4").Value = 11   # was 10
$ws.Range("L7").Value = 152   # was 151

# --- Avalon Park ---
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L6").Value = 30   # was 29
$ws.Range("L7").Value = 124   # was 123

# --- Portage Park ---
$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L2").Value = 62   # was 60
$ws.Range("L7").Value = 171   # was 169

# --- Streeterville ---
$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L6").Value = 18   # was 17
$ws.Range("L7").Value = 136   # was 135

# --- Little Italy, UIC ---
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L6").Value = 56   # was 55
$ws.Range("L7").Value = 266   # was 265

# --- Mckinley Park ---
$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("L2").Value = 25   # was 24
$ws.Range("L7").Value = 73   # was 72

# --- Oakland ---
$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("L3").Value = 19   # was 18
$ws.Range("L7").Value = 55   # was 54

# --- Riverdale ---
$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L3").Value = 47   # was 46
$ws.Range("L7").Value = 144   # was 143

# --- Andersonville ---
$ws = $wb.Worksheets.Item("Andersonville")
$ws.Range("L6").Value = 3   # was 2
$ws.Range("L7").Value = 14   # was 13

# --- Ukrainian Village ---
$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("L3").Value = 15   # was 14
$ws.Range("L4").Value = 9   # was 8
$ws.Range("L7").Value = 59   # was 57

# --- Millenium Park ---
$ws = $wb.Worksheets.Item("Millenium Park")
$ws.Range("L6").Value = 5   # was 6
$ws.Range("L7").Value = 11   # was 12
